$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 83 (id 81) and Row 84 (id 82) have their data (columns B, F:AC) swapped.
# Column A (id) and columns C, D, E (league/date, identical on both rows) stay put.

$row1 = 83
$row2 = 84

# Columns to swap: B, then F through AC (skip C, D, E)
$cols = @("B")
$cols += @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

foreach ($col in $cols) {
    $addr1 = "$col$row1"
    $addr2 = "$col$row2"
    $v1 = $ws.Range($addr1).Value2
    $v2 = $ws.Range($addr2).Value2
    $ws.Range($addr1).Value = $v2
    $ws.Range($addr2).Value = $v1
}

Write-Output "swap complete"
